$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 580, shifting existing rows 580:621 down to 581:622
$ws.Rows(580).Insert()

# Populate the newly inserted row with the new record.
# Use an apostrophe prefix so the date-like text stays plain text (matches
# the rest of the column, which stores dates as text, not date serials),
# then strip the resulting cell formatting so no stray style index is left
# behind (ClearFormats keeps the cell's stored value/type, only formatting).
$ws.Range("A580").Value = "'2026/01/07"
$ws.Range("A580").ClearFormats()
$ws.Range("B580").Value = "水"
$ws.Range("C580").Value = 3
$ws.Range("D580").Value = 201
